$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition listing)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 511
$ws1.Range("F3").Value = 3425

# Sheet "全部类型" (all types listing)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 511
$ws4.Range("F3").Value = 3426
